# Rows 5, 7 and 8 of the "Artfynd" sheet got their observation data rotated:
#   new row 5 <- old row 8 data
#   new row 7 <- old row 5 data
#   new row 8 <- old row 7 data
# (row 6 is untouched). Because every destination value is fully known up
# front, we just write the final literal values into each cell instead of
# doing a positional swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 5 (becomes the former row-8 record: Flagellkvastmossa / Dicranum
# flagellare, with the "med groddkorn" age/stage note)
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 111934066
$ws.Range("B5").Value = 93289
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 2170
$ws.Range("F5").Value = "Flagellkvastmossa"
$ws.Range("G5").Value = "Dicranum flagellare"
$ws.Range("H5").Value = "Hedw."

# I5 / J5 had "1" / "fruktkroppar" - they become blank text cells
$ws.Range("I5").Value = "'"
$ws.Range("I5").Style = "Normal"
$ws.Range("J5").Value = "'"
$ws.Range("J5").Style = "Normal"

# K5 was blank - gains the age/stage note
$ws.Range("K5").Value = "med groddkorn"

# L5 was absent/blank - becomes a present-but-empty text cell
$ws.Range("L5").Value = "'"
$ws.Range("L5").Style = "Normal"

$ws.Range("P5").Value = "Skogen N om Dye, I2-Skogen, Vrm"
$ws.Range("Q5").Value = 413590.3038565172
$ws.Range("R5").Value = 6586912.201658082

# ---------------------------------------------------------------------
# Row 7 (becomes the former row-5 record: Motaggsvamp / Sarcodon
# squamosus, reported as 1 fruktkropp)
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 111934086
$ws.Range("B7").Value = 90689
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 5966
$ws.Range("F7").Value = "Motaggsvamp"
$ws.Range("G7").Value = "Sarcodon squamosus"
$ws.Range("H7").Value = "(Schaeff.) Quél."

# I7 / J7 were blank text cells - gain "1" / "fruktkroppar"
# (force text storage for "1" with a leading apostrophe so it round-trips
# as a text cell rather than a number, matching the original sheet)
$ws.Range("I7").Value = "'1"
$ws.Range("I7").Style = "Normal"
$ws.Range("J7").Value = "fruktkroppar"

# K7 had "med groddkorn" - becomes a blank text cell
$ws.Range("K7").Value = "'"
$ws.Range("K7").Style = "Normal"

# L7 was a present-but-empty text cell - becomes fully absent/blank
$ws.Range("L7").Value = ""

$ws.Range("P7").Value = "Tallskogen N om Dye, I2-Skogen, Vrm"
$ws.Range("Q7").Value = 413681.2082122188
$ws.Range("R7").Value = 6586805.223123537

# AC7 held a public comment - it disappears entirely
$ws.Range("AC7").Value = ""

# ---------------------------------------------------------------------
# Row 8 (becomes the former row-7 record: Flagellkvastmossa / Dicranum
# flagellare, with the public comment about being common along the trail)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = 111934059
$ws.Range("Q8").Value = 413639.6308819132
$ws.Range("R8").Value = 6586793.951973591

# AC8 was absent/blank - gains the public comment
$ws.Range("AC8").Value = "Rätt riklig längs stigen"
